$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-03-11 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-03-12 Sunday", 2) | Out-Null
$d.Content.Find.Execute("98-61=37", $true, $false, $false, $false, $false, $true, 1, $false, "98-4=94", 2) | Out-Null
$d.Content.Find.Execute("77-50=27", $true, $false, $false, $false, $false, $true, 1, $false, "57+42=99", 2) | Out-Null
$d.Content.Find.Execute("78-21=57", $true, $false, $false, $false, $false, $true, 1, $false, "16-4=12", 2) | Out-Null
$d.Content.Find.Execute("39+0=39", $true, $false, $false, $false, $false, $true, 1, $false, "86-73=13", 2) | Out-Null
$d.Content.Find.Execute("64-24=40", $true, $false, $false, $false, $false, $true, 1, $false, "82-26=56", 2) | Out-Null
$d.Content.Find.Execute("11+66=77", $true, $false, $false, $false, $false, $true, 1, $false, "8+30=38", 2) | Out-Null
$d.Content.Find.Execute("83-79=4", $true, $false, $false, $false, $false, $true, 1, $false, "51+32=83", 2) | Out-Null
$d.Content.Find.Execute("0+63=63", $true, $false, $false, $false, $false, $true, 1, $false, "24+56=80", 2) | Out-Null
$d.Content.Find.Execute("55-49=6", $true, $false, $false, $false, $false, $true, 1, $false, "62-27=35", 2) | Out-Null
$d.Content.Find.Execute("64-21=43", $true, $false, $false, $false, $false, $true, 1, $false, "89-68=21", 2) | Out-Null
$d.Content.Find.Execute("43-18=25", $true, $false, $false, $false, $false, $true, 1, $false, "68-53=15", 2) | Out-Null
$d.Content.Find.Execute("93-15=78", $true, $false, $false, $false, $false, $true, 1, $false, "18-8=10", 2) | Out-Null
$d.Content.Find.Execute("82-37=45", $true, $false, $false, $false, $false, $true, 1, $false, "5+18=23", 2) | Out-Null
$d.Content.Find.Execute("31+22=53", $true, $false, $false, $false, $false, $true, 1, $false, "45-45=0", 2) | Out-Null
$d.Content.Find.Execute("98-51=47", $true, $false, $false, $false, $false, $true, 1, $false, "81-3=78", 2) | Out-Null
$d.Content.Find.Execute("99-50=49", $true, $false, $false, $false, $false, $true, 1, $false, "5+61=66", 2) | Out-Null
$d.Content.Find.Execute("26+16=42", $true, $false, $false, $false, $false, $true, 1, $false, "70-28=42", 2) | Out-Null
$d.Content.Find.Execute("58+20=78", $true, $false, $false, $false, $false, $true, 1, $false, "72+2=74", 2) | Out-Null
$d.Content.Find.Execute("26+39=65", $true, $false, $false, $false, $false, $true, 1, $false, "4+12=16", 2) | Out-Null
$d.Content.Find.Execute("19+68=87", $true, $false, $false, $false, $false, $true, 1, $false, "57-4=53", 2) | Out-Null
$d.Content.Find.Execute("80-59=21", $true, $false, $false, $false, $false, $true, 1, $false, "61-8=53", 2) | Out-Null
$d.Content.Find.Execute("99-11=88", $true, $false, $false, $false, $false, $true, 1, $false, "66+14=80", 2) | Out-Null
$d.Content.Find.Execute("13-2=11", $true, $false, $false, $false, $false, $true, 1, $false, "89-32=57", 2) | Out-Null
$d.Content.Find.Execute("61+37=98", $true, $false, $false, $false, $false, $true, 1, $false, "3-2=1", 2) | Out-Null
$d.Content.Find.Execute("97-41=56", $true, $false, $false, $false, $false, $true, 1, $false, "74-9=65", 2) | Out-Null
$d.Content.Find.Execute("65-63=2", $true, $false, $false, $false, $false, $true, 1, $false, "48+23=71", 2) | Out-Null
$d.Content.Find.Execute("45+17=62", $true, $false, $false, $false, $false, $true, 1, $false, "44+39=83", 2) | Out-Null
$d.Content.Find.Execute("7+62=69", $true, $false, $false, $false, $false, $true, 1, $false, "19-0=19", 2) | Out-Null
$d.Content.Find.Execute("9+29=38", $true, $false, $false, $false, $false, $true, 1, $false, "65-34=31", 2) | Out-Null
$d.Content.Find.Execute("61+32=93", $true, $false, $false, $false, $false, $true, 1, $false, "81+9=90", 2) | Out-Null
$d.Content.Find.Execute("30-30=0", $true, $false, $false, $false, $false, $true, 1, $false, "56-13=43", 2) | Out-Null
$d.Content.Find.Execute("18+45=63", $true, $false, $false, $false, $false, $true, 1, $false, "16+39=55", 2) | Out-Null
$d.Content.Find.Execute("26-10=16", $true, $false, $false, $false, $false, $true, 1, $false, "66+31=97", 2) | Out-Null
$d.Content.Find.Execute("57+7=64", $true, $false, $false, $false, $false, $true, 1, $false, "78+14=92", 2) | Out-Null
$d.Content.Find.Execute("99-73=26", $true, $false, $false, $false, $false, $true, 1, $false, "95-83=12", 2) | Out-Null
$d.Content.Find.Execute("86-9=77", $true, $false, $false, $false, $false, $true, 1, $false, "54+0=54", 2) | Out-Null
$d.Content.Find.Execute("6+52=58", $true, $false, $false, $false, $false, $true, 1, $false, "71-6=65", 2) | Out-Null
$d.Content.Find.Execute("82+1=83", $true, $false, $false, $false, $false, $true, 1, $false, "4+24=28", 2) | Out-Null
$d.Content.Find.Execute("48+9=57", $true, $false, $false, $false, $false, $true, 1, $false, "47-18=29", 2) | Out-Null
$d.Content.Find.Execute("62+2=64", $true, $false, $false, $false, $false, $true, 1, $false, "30-22=8", 2) | Out-Null
$d.Content.Find.Execute("98-29=69", $true, $false, $false, $false, $false, $true, 1, $false, "78+13=91", 2) | Out-Null
$d.Content.Find.Execute("95-47=48", $true, $false, $false, $false, $false, $true, 1, $false, "34+7=41", 2) | Out-Null
$d.Content.Find.Execute("8+19=27", $true, $false, $false, $false, $false, $true, 1, $false, "30-20=10", 2) | Out-Null
$d.Content.Find.Execute("62+23=85", $true, $false, $false, $false, $false, $true, 1, $false, "14+32=46", 2) | Out-Null
$d.Content.Find.Execute("25-2=23", $true, $false, $false, $false, $false, $true, 1, $false, "12+21=33", 2) | Out-Null
$d.Content.Find.Execute("29+37=66", $true, $false, $false, $false, $false, $true, 1, $false, "60-2=58", 2) | Out-Null
$d.Content.Find.Execute("84+14=98", $true, $false, $false, $false, $false, $true, 1, $false, "58+33=91", 2) | Out-Null
$d.Content.Find.Execute("80-75=5", $true, $false, $false, $false, $false, $true, 1, $false, "61+31=92", 2) | Out-Null
$d.Content.Find.Execute("87-12=75", $true, $false, $false, $false, $false, $true, 1, $false, "78-13=65", 2) | Out-Null
$d.Content.Find.Execute("40-0=40", $true, $false, $false, $false, $false, $true, 1, $false, "73-66=7", 2) | Out-Null
$d.Content.Find.Execute("54+16=70", $true, $false, $false, $false, $false, $true, 1, $false, "78-73=5", 2) | Out-Null
$d.Content.Find.Execute("30+40=70", $true, $false, $false, $false, $false, $true, 1, $false, "71+1=72", 2) | Out-Null
$d.Content.Find.Execute("38+57=95", $true, $false, $false, $false, $false, $true, 1, $false, "41+10=51", 2) | Out-Null
$d.Content.Find.Execute("35+13=48", $true, $false, $false, $false, $false, $true, 1, $false, "59-20=39", 2) | Out-Null
$d.Content.Find.Execute("35-7=28", $true, $false, $false, $false, $false, $true, 1, $false, "91-38=53", 2) | Out-Null
$d.Content.Find.Execute("9+63=72", $true, $false, $false, $false, $false, $true, 1, $false, "30+63=93", 2) | Out-Null
$d.Content.Find.Execute("20+3=23", $true, $false, $false, $false, $false, $true, 1, $false, "80-66=14", 2) | Out-Null
$d.Content.Find.Execute("82+11=93", $true, $false, $false, $false, $false, $true, 1, $false, "56+17=73", 2) | Out-Null
$d.Content.Find.Execute("70-50=20", $true, $false, $false, $false, $false, $true, 1, $false, "44-27=17", 2) | Out-Null
$d.Content.Find.Execute("47-23=24", $true, $false, $false, $false, $false, $true, 1, $false, "66-16=50", 2) | Out-Null
$d.Content.Find.Execute("89-65=24", $true, $false, $false, $false, $false, $true, 1, $false, "73-44=29", 2) | Out-Null
$d.Content.Find.Execute("98-93=5", $true, $false, $false, $false, $false, $true, 1, $false, "2+85=87", 2) | Out-Null
$d.Content.Find.Execute("1+57=58", $true, $false, $false, $false, $false, $true, 1, $false, "59+23=82", 2) | Out-Null
$d.Content.Find.Execute("25+38=63", $true, $false, $false, $false, $false, $true, 1, $false, "94-84=10", 2) | Out-Null
$d.Content.Find.Execute("32+55=87", $true, $false, $false, $false, $false, $true, 1, $false, "50-3=47", 2) | Out-Null
$d.Content.Find.Execute("23-6=17", $true, $false, $false, $false, $false, $true, 1, $false, "98-34=64", 2) | Out-Null
$d.Content.Find.Execute("94-80=14", $true, $false, $false, $false, $false, $true, 1, $false, "90-88=2", 2) | Out-Null
$d.Content.Find.Execute("88+6=94", $true, $false, $false, $false, $false, $true, 1, $false, "40-11=29", 2) | Out-Null
$d.Content.Find.Execute("52+30=82", $true, $false, $false, $false, $false, $true, 1, $false, "13+53=66", 2) | Out-Null
$d.Content.Find.Execute("91-16=75", $true, $false, $false, $false, $false, $true, 1, $false, "85-41=44", 2) | Out-Null
$d.Content.Find.Execute("92-34=58", $true, $false, $false, $false, $false, $true, 1, $false, "37+23=60", 2) | Out-Null
$d.Content.Find.Execute("25+23=48", $true, $false, $false, $false, $false, $true, 1, $false, "36+42=78", 2) | Out-Null
$d.Content.Find.Execute("19+29=48", $true, $false, $false, $false, $false, $true, 1, $false, "85+3=88", 2) | Out-Null
$d.Content.Find.Execute("13-7=6", $true, $false, $false, $false, $false, $true, 1, $false, "24+28=52", 2) | Out-Null
$d.Content.Find.Execute("64-1=63", $true, $false, $false, $false, $false, $true, 1, $false, "20+27=47", 2) | Out-Null
$d.Content.Find.Execute("76-26=50", $true, $false, $false, $false, $false, $true, 1, $false, "15+37=52", 2) | Out-Null
$d.Content.Find.Execute("93-32=61", $true, $false, $false, $false, $false, $true, 1, $false, "53+15=68", 2) | Out-Null
$d.Content.Find.Execute("24+13=37", $true, $false, $false, $false, $false, $true, 1, $false, "58-49=9", 2) | Out-Null
$d.Content.Find.Execute("84-70=14", $true, $false, $false, $false, $false, $true, 1, $false, "13+64=77", 2) | Out-Null
$d.Content.Find.Execute("57-19=38", $true, $false, $false, $false, $false, $true, 1, $false, "89-8=81", 2) | Out-Null
$d.Content.Find.Execute("42+48=90", $true, $false, $false, $false, $false, $true, 1, $false, "36+45=81", 2) | Out-Null
$d.Content.Find.Execute("83-12=71", $true, $false, $false, $false, $false, $true, 1, $false, "11+21=32", 2) | Out-Null
$d.Content.Find.Execute("77+5=82", $true, $false, $false, $false, $false, $true, 1, $false, "33+30=63", 2) | Out-Null
$d.Content.Find.Execute("70-2=68", $true, $false, $false, $false, $false, $true, 1, $false, "84-67=17", 2) | Out-Null
$d.Content.Find.Execute("78-36=42", $true, $false, $false, $false, $false, $true, 1, $false, "78+2=80", 2) | Out-Null
$d.Content.Find.Execute("25-9=16", $true, $false, $false, $false, $false, $true, 1, $false, "8+43=51", 2) | Out-Null
$d.Content.Find.Execute("84-59=25", $true, $false, $false, $false, $false, $true, 1, $false, "53+15=68", 2) | Out-Null
$d.Content.Find.Execute("85-3=82", $true, $false, $false, $false, $false, $true, 1, $false, "52-46=6", 2) | Out-Null
$d.Content.Find.Execute("33+29=62", $true, $false, $false, $false, $false, $true, 1, $false, "5+38=43", 2) | Out-Null
$d.Content.Find.Execute("93-18=75", $true, $false, $false, $false, $false, $true, 1, $false, "27+64=91", 2) | Out-Null
$d.Content.Find.Execute("40-35=5", $true, $false, $false, $false, $false, $true, 1, $false, "37-13=24", 2) | Out-Null
$d.Content.Find.Execute("58-3=55", $true, $false, $false, $false, $false, $true, 1, $false, "27+11=38", 2) | Out-Null
$d.Content.Find.Execute("74-22=52", $true, $false, $false, $false, $false, $true, 1, $false, "83-15=68", 2) | Out-Null
$d.Content.Find.Execute("21+68=89", $true, $false, $false, $false, $false, $true, 1, $false, "31+18=49", 2) | Out-Null
$d.Content.Find.Execute("94-38=56", $true, $false, $false, $false, $false, $true, 1, $false, "5+93=98", 2) | Out-Null
$d.Content.Find.Execute("49-16=33", $true, $false, $false, $false, $false, $true, 1, $false, "42-5=37", 2) | Out-Null
$d.Content.Find.Execute("60-46=14", $true, $false, $false, $false, $false, $true, 1, $false, "19+19=38", 2) | Out-Null
$d.Content.Find.Execute("6+26=32", $true, $false, $false, $false, $false, $true, 1, $false, "87-50=37", 2) | Out-Null
$d.Content.Find.Execute("92-29=63", $true, $false, $false, $false, $false, $true, 1, $false, "38-17=21", 2) | Out-Null
$d.Content.Find.Execute("96-68=28", $true, $false, $false, $false, $false, $true, 1, $false, "33+43=76", 2) | Out-Null
